$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.339434623718262
$ws.Range("B1").Value = 2.466507911682129
$ws.Range("C1").Value = 4.842360973358154
$ws.Range("D1").Value = 2.441564321517944
$ws.Range("E1").Value = 0.9265078902244568
